# Recompute "Anteil Blähton" / "Anteil Erde" as eighths-of-40 fractions
# instead of quarters (8/40 = 0.2, 32/40 = 0.8) - this feeds the rest of
# the dependent calculation chain (B9, B10, B15, B17, B18) automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Formula = "=8/40"
$ws.Range("B13").Formula = "=32/40"

# Drop the stray "Materialien" header row (row 23) that was left dangling
# at the bottom of the sheet, and fix up the "Gesamtvolumen Blähton" label
# that had been (re)written further down the sheet.
$ws.Range("A23").ClearContents()
$ws.Range("A10").Value = "Gesamtvolumen Blähton"

# Update the active selection to reflect where editing left off.
$ws.Range("E6:E7").Select() | Out-Null
